$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64, shifting existing rows 64..150 down to 65..151
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record
$ws.Cells.Item(64, 1).Value = 10
$ws.Cells.Item(64, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value = "La Araucanía"
$ws.Cells.Item(64, 4).Value = 44413
$ws.Cells.Item(64, 5).Value = 9
$ws.Cells.Item(64, 6).Value = 100112001
$ws.Cells.Item(64, 7).Value = "Berenjena"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 200
$ws.Cells.Item(64, 11).Value = 13000
$ws.Cells.Item(64, 12).Value = 14000
$ws.Cells.Item(64, 13).Value = 13500
$ws.Cells.Item(64, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 225
$ws.Cells.Item(64, 17).Value = 60
$ws.Cells.Item(64, 18).Value = "Hortaliza"
